$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Row 17 / 18 coin swap (WrappedBTC <-> WrappedEther) ---
Set-TextValue $ws.Range("B17") 'WrappedEther'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D17") '3.649.19'
Set-TextValue $ws.Range("E17") '  +7.36%  '

Set-TextValue $ws.Range("B18") 'WrappedBTC'
Set-TextValue $ws.Range("C18") 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range("D18") '71.949.39'
Set-TextValue $ws.Range("E18") '  +2.59%  '

# --- Remaining per-cell value updates ---
Set-TextValue $ws.Range("D2") '71.884.44'
Set-TextValue $ws.Range("E2") '  +2.49%  '
Set-TextValue $ws.Range("D3") '3.658.06'
Set-TextValue $ws.Range("E3") '  +7.43%  '
Set-TextValue $ws.Range("E4") '  -0.02%  '
Set-TextValue $ws.Range("D5") '592.02'
Set-TextValue $ws.Range("E5") '  +1.35%  '
Set-TextValue $ws.Range("D6") '181.33'
Set-TextValue $ws.Range("E6") '  -0.42%  '
Set-TextValue $ws.Range("D7") '3.650.42'
Set-TextValue $ws.Range("E7") '  +7.38%  '
Set-TextValue $ws.Range("E8") '  +3.55%  '
Set-TextValue $ws.Range("E9") '  +0.02%  '
Set-TextValue $ws.Range("E10") '  +0.80%  '
Set-TextValue $ws.Range("D11") '0.611'
Set-TextValue $ws.Range("E11") '  +2.73%  '
Set-TextValue $ws.Range("D12") '49.97'
Set-TextValue $ws.Range("E12") '  +2.23%  '
Set-TextValue $ws.Range("E13") '  -0.36%  '
Set-TextValue $ws.Range("D14") '682.75'
Set-TextValue $ws.Range("E14") '  -0.93%  '
Set-TextValue $ws.Range("D15") '4.236.46'
Set-TextValue $ws.Range("E15") '  +7.04%  '
Set-TextValue $ws.Range("D16") '9.08'
Set-TextValue $ws.Range("E16") '  +4.21%  '
Set-TextValue $ws.Range("E19") '  +1.70%  '
Set-TextValue $ws.Range("E20") '  +2.98%  '
Set-TextValue $ws.Range("E21") '  +2.35%  '
Set-TextValue $ws.Range("D22") '0.944'
Set-TextValue $ws.Range("E22") '  +2.77%  '
Set-TextValue $ws.Range("D23") '6.14'
Set-TextValue $ws.Range("E23") '  +14.28%  '
Set-TextValue $ws.Range("D24") '18.05'
Set-TextValue $ws.Range("E24") '  +3.66%  '
Set-TextValue $ws.Range("D25") '103.56'
Set-TextValue $ws.Range("E25") '  +0.90%  '
Set-TextValue $ws.Range("E26") '  +2.73%  '
Set-TextValue $ws.Range("E27") '  +4.15%  '
Set-TextValue $ws.Range("D28") '10.21'
Set-TextValue $ws.Range("E28") '  +4.99%  '
Set-TextValue $ws.Range("D29") '35.41'
Set-TextValue $ws.Range("E29") '  +4.39%  '
Set-TextValue $ws.Range("E30") '  +4.39%  '
Set-TextValue $ws.Range("D31") '7.35'
Set-TextValue $ws.Range("E31") '  +5.09%  '
Set-TextValue $ws.Range("E32") '  +12.53%  '
Set-TextValue $ws.Range("D33") '580.48'
Set-TextValue $ws.Range("E33") '  +4.11%  '
Set-TextValue $ws.Range("D34") '11.35'
Set-TextValue $ws.Range("E34") '  +1.73%  '
Set-TextValue $ws.Range("E35") '  +1.89%  '
Set-TextValue $ws.Range("D36") '59.56'
Set-TextValue $ws.Range("E36") '  +1.46%  '
Set-TextValue $ws.Range("E37") '  -0.05%  '
Set-TextValue $ws.Range("D38") '3.731.67'
Set-TextValue $ws.Range("E38") '  +1.89%  '
Set-TextValue $ws.Range("E39") '  +2.94%  '
Set-TextValue $ws.Range("D40") '35.75'
Set-TextValue $ws.Range("E40") '  +0.16%  '
Set-TextValue $ws.Range("D41") '0.0₃0763'
Set-TextValue $ws.Range("E41") '  +3.61%  '
Set-TextValue $ws.Range("D42") '3.48'
Set-TextValue $ws.Range("E42") '  +4.22%  '
Set-TextValue $ws.Range("D43") '0.0468'
Set-TextValue $ws.Range("E43") '  +8.78%  '
Set-TextValue $ws.Range("D44") '2.81'
Set-TextValue $ws.Range("E44") '  +2.05%  '
Set-TextValue $ws.Range("D45") '0.348'
Set-TextValue $ws.Range("E45") '  +2.22%  '
Set-TextValue $ws.Range("D46") '3.38'
Set-TextValue $ws.Range("E46") '  +0.59%  '
Set-TextValue $ws.Range("E47") '  +5.19%  '
Set-TextValue $ws.Range("E48") '  +3.05%  '
Set-TextValue $ws.Range("E49") '  +3.99%  '
Set-TextValue $ws.Range("E50") '  -0.03%  '
Set-TextValue $ws.Range("D51") '132.73'
Set-TextValue $ws.Range("E51") '  +1.61%  '
